# BKQBF_TAB-5 PackingList update:
# Insert a new "alternate code" row right after the existing
# "GARFO SOBREMESA - 1000UN" row (row 71), duplicating that product with a
# new item code (31543) and flagging it with the sheet's existing
# "new code / highlighted" row style (yellow fill), same as the other
# special rows already present further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 72 (pushes old row 72.. down to 73..)
$ws.Rows.Item(72).Insert()

# Copy the formatting (only, not values) of the existing "highlighted /
# alternate code" row (the PAO BK row, now at 123 after the insert above)
# onto the new row A72:E72 so it gets the same yellow-highlight style used
# elsewhere in this sheet for this kind of entry.
$ws.Range("A123:E123").Copy()
$ws.Range("A72:E72").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's data: same product/description and quantities as
# the row above it (GARFO SOBREMESA - 1000UN), but under the new code.
$ws.Range("A72").Value2 = 31543
$ws.Range("B72").Value2 = "GARFO SOBREMESA - 1000UN"
$ws.Range("C72").Value2 = 1000
$ws.Range("D72").Value2 = 50
$ws.Range("E72").Value2 = 1

# Update the sheet's active selection to match the saved view.
$ws.Range("F16").Select()
